# Workbook audit cleanup:
#  - "optimization_parameters" had a stray leftover row (a "Sheet" label
#    with orphaned numbers 3 / 4 next to it) that doesn't belong in the
#    parameter table; delete it so the sheet's data block closes up.
#  - Leave the cursor on a few sheets as they were reviewed, finishing on
#    "threshold_b" as the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Glance at network_weights while auditing.
$wsWeights = $wb.Worksheets.Item("network_weights")
$wsWeights.Activate() | Out-Null
$wsWeights.Range("E9").Select() | Out-Null

# The actual fix: remove the stray "Sheet" row from optimization_parameters.
$wsParams = $wb.Worksheets.Item("optimization_parameters")
$wsParams.Activate() | Out-Null
$wsParams.Rows("16:16").Delete() | Out-Null
$wsParams.Range("A16:XFD16").Select() | Out-Null

# Finish the audit on threshold_b, which becomes the active sheet.
$wsThreshold = $wb.Worksheets.Item("threshold_b")
$wsThreshold.Activate() | Out-Null
$wsThreshold.Range("B42").Select() | Out-Null
